$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.24056878790175773
$ws.Range("A2").Value = -0.0059999999846880314
$ws.Range("A3").Value = -0.0039999999854689605
$ws.Range("A4").Value = -0.0079999999742987882
$ws.Range("A5").Value = -0.0029999999848913106
$ws.Range("A6").Value = -0.0019999999836937121
$ws.Range("A7").Value = -0.0099999999643656245
$ws.Range("A8").Value = 0.006905092596982243
$ws.Range("A9").Value = -0.0019999999842168492
$ws.Range("A10").Value = -0.0019999999844859673
$ws.Range("A11").Value = -0.0029999999821077594
$ws.Range("A12").Value = -0.0035050534445764114
$ws.Range("A13").Value = -0.003499999981841917
$ws.Range("A14").Value = -0.0079999999713447068
$ws.Range("A15").Value = -0.00099999998889988007
$ws.Range("A16").Value = -0.0019999999873157037
$ws.Range("A17").Value = 0.034479060806861916
$ws.Range("A18").Value = -0.0039999999839865907
$ws.Range("A19").Value = -0.0039999999894995142
$ws.Range("A20").Value = -0.0039999999888049587
$ws.Range("A21").Value = -0.0039999999886841664
$ws.Range("A22").Value = -0.0039999999886077831
$ws.Range("A23").Value = -0.061931894804417631
$ws.Range("A24").Value = -0.019999999943599356
$ws.Range("A25").Value = -0.01999999994289503
$ws.Range("A26").Value = -0.0024999999817811869
$ws.Range("A27").Value = -0.002499999980810852
$ws.Range("A28").Value = -0.0019999999779578559
$ws.Range("A29").Value = -0.006999999962643777
$ws.Range("A30").Value = -0.059999999832742734
$ws.Range("A31").Value = -0.0069999999599410501
$ws.Range("A32").Value = -0.0099999999525177685
$ws.Range("A33").Value = -0.0039999999669504405
